$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.2303363333333333
$ws.Range("N2").Value = 0.691009
$ws.Range("O2").Value = 0.0420565315194687
$ws.Range("P2").Value = 0.0420565315194687
$ws.Range("Q2").Value = 0.4250969128682223
$ws.Range("R2").Value = 3.825872215814
$ws.Range("S2").Value = 0.0420565315194687
$ws.Range("T2").Value = 0.0420565315194687

# Row 3 updates
$ws.Range("O3").Value = 0.8440851393264226
$ws.Range("P3").Value = 0.8440851393264227
$ws.Range("S3").Value = 0.8440851393264226
$ws.Range("T3").Value = 0.8440851393264227

# Row 4 updates
$ws.Range("M4").Value = 0.6235823333333333
$ws.Range("N4").Value = 1.870747
$ws.Range("O4").Value = 0.1138583291541087
$ws.Range("P4").Value = 0.1138583291541087
$ws.Range("Q4").Value = 1.150851543840222
$ws.Range("R4").Value = 10.357663894562
$ws.Range("S4").Value = 0.1138583291541087
$ws.Range("T4").Value = 0.1138583291541087
